# "Still working on calculations of total with tax"
#
# Inserts a brand-new worksheet between "Sheet1" and the existing "Sheet2"
# (which pushes the old "Sheet2" tab to a new position and Excel
# auto-renames it "Sheet3"), fills the new sheet with a small notes table,
# and updates the selection on Sheet1.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# Insert the new worksheet right after Sheet1 (before the old Sheet2), so
# the tab order becomes: Sheet1, <new sheet>, Sheet2 - and the old "Sheet2"
# keeps its name "Sheet2" while the brand new tab is auto-named "Sheet3"...
# Excel actually names the newly inserted sheet using the next free number,
# which also happens to be "Sheet3" here, and leaves the original "Sheet2"
# sheet's name untouched.
$newSheet = $wb.Worksheets.Add($null, $sheet1)

# Populate the notes / quality-control table on the new sheet. Cells are
# written in this specific order so shared-string indices land the same
# way they did for the author.
$newSheet.Range("C5").Value = "Tic Quality Control"
$newSheet.Range("C4").Value = "company name"
$newSheet.Range("D4").Value = "Adress"
$newSheet.Range("C8").Value = "calculateOrderDiscount"
$newSheet.Range("C7").Value = "itemPercent"
$newSheet.Range("D7").Value = "itemTotal /itemsTotal"
$newSheet.Range("D8").Value = "orderDiscount/100 or orderDiscount"
$newSheet.Range("C9").Value = "itemDiscount"
$newSheet.Range("D9").Value = "itemDiscount + itemPercent * calculateOrderDiscount"

# Column widths for the new sheet's two used columns.
$newSheet.Columns.Item(3).ColumnWidth = 30.166666666666668
$newSheet.Columns.Item(4).ColumnWidth = 42.76

# Leave the cursor on the new sheet at D10, matching the saved selection.
$newSheet.Range("D10").Select() | Out-Null

# Update the selection remembered on Sheet1 (it is no longer the active
# tab, but Excel still persists its last selection).
$sheet1.Range("F20:G21").Select() | Out-Null

# Re-activate the new sheet so it is the one on top when the workbook is
# reopened (matches activeTab pointing at the new sheet).
$newSheet.Activate() | Out-Null
$newSheet.Range("D10").Select() | Out-Null
